{"js": "// Load all paragraphs in the document body so we can find the two\n// placeholder paragraphs (\"<<judgeRecital>>\" and \"<<generalOrder>>\") by\n// their text content (more robust than relying on fixed indices).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet judgeRecitalPara = null;\nlet generalOrderPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text.trim();\n  if (t === \"<<judgeRecital>>\") {\n    judgeRecitalPara = p;\n  } else if (t === \"<<generalOrder>>\") {\n    generalOrderPara = p;\n  }\n}\n\nif (!judgeRecitalPara) {\n  throw new Error(\"Could not find the '<<judgeRecital>>' placeholder paragraph.\");\n}\nif (!generalOrderPara) {\n  throw new Error(\"Could not find the '<<generalOrder>>' placeholder paragraph.\");\n}\n\n// 1) Turn the \"<<judgeRecital>>\" paragraph into the new recital sentence\n//    that references the applicant name and application date.\njudgeRecitalPara.insertText(\n  \"Upon the application of <<applicantName>> dated <<applicationDate>> and upon considering the information provided by the parties:\",\n  \"Replace\"\n);\n\n// 2) Leave the \"<<generalOrder>>\" paragraph itself untouched, but append a\n//    new \"Reasons for decision:\" block (with surrounding blank paragraphs)\n//    followed by a new \"<<reasonForDecision>>\" placeholder paragraph right\n//    after it.\nconst blank1 = generalOrderPara.insertParagraph(\"\", \"After\");\nconst reasonsHeading = blank1.insertParagraph(\"Reasons for decision:\", \"After\");\nconst blank2 = reasonsHeading.insertParagraph(\"\", \"After\");\nblank2.insertParagraph(\"<<reasonForDecision>>\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Turn the \"<<judgeRecital>>\" paragraph into the new recital sentence\n#    that references the applicant name and application date. A plain\n#    Find & Replace keeps the surrounding paragraph mark intact.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$newRecital = \"Upon the application of <<applicantName>> dated <<applicationDate>> and upon considering the information provided by the parties:\"\n$find.Text = \"<<judgeRecital>>\"\n$find.Replacement.Text = $newRecital\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $found) {\n    throw \"Could not find the '<<judgeRecital>>' placeholder paragraph.\"\n}\n\n# Locate the (untouched) \"<<generalOrder>>\" placeholder paragraph so new\n# content can be appended right after it.\n$generalOrderPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\")\n    if ($t -eq \"<<generalOrder>>\") {\n        $generalOrderPara = $p\n        break\n    }\n}\nif ($null -eq $generalOrderPara) {\n    throw \"Could not find the '<<generalOrder>>' placeholder paragraph.\"\n}\n\n# 2) Leave the \"<<generalOrder>>\" paragraph itself untouched, but append a\n#    new \"Reasons for decision:\" block (with surrounding blank paragraphs)\n#    followed by a new \"<<reasonForDecision>>\" placeholder paragraph right\n#    after it.\n$generalOrderPara.Range.InsertParagraphAfter()\n$blank1 = $generalOrderPara.Next()\n\n$blank1.Range.InsertParagraphAfter()\n$reasonsHeading = $blank1.Next()\n$reasonsHeading.Range.Text = \"Reasons for decision:\"\n\n$reasonsHeading.Range.InsertParagraphAfter()\n$blank2 = $reasonsHeading.Next()\n\n$blank2.Range.InsertParagraphAfter()\n$reasonForDecisionPara = $blank2.Next()\n$reasonForDecisionPara.Range.Text = \"<<reasonForDecision>>\"\n"}
